$d = $word.ActiveDocument

$replacements = @(
    @{ old = "66×25="; new = "99×26=" },
    @{ old = "98×27="; new = "39×89=" },
    @{ old = "76×27="; new = "71×82=" },
    @{ old = "41×29="; new = "71×36=" },
    @{ old = "70×12="; new = "27×12=" },
    @{ old = "34×54="; new = "90×22=" },
    @{ old = "98×52="; new = "15×82=" },
    @{ old = "75×34="; new = "19×43=" },
    @{ old = "44×71="; new = "84×62=" },
    @{ old = "80×67="; new = "34×20=" },
    @{ old = "29×18="; new = "22×63=" },
    @{ old = "96×62="; new = "41×20=" },
    @{ old = "30×51="; new = "66×37=" },
    @{ old = "75×45="; new = "35×22=" },
    @{ old = "66×95="; new = "87×60=" },
    @{ old = "57×17="; new = "39×22=" },
    @{ old = "16×91="; new = "44×90=" },
    @{ old = "88×81="; new = "53×11=" },
    @{ old = "13×39="; new = "26×64=" },
    @{ old = "68×16="; new = "11×57=" },
    @{ old = "30×64="; new = "21×85=" },
    @{ old = "99×43="; new = "86×27=" },
    @{ old = "66×32="; new = "50×81=" },
    @{ old = "46×70="; new = "42×35=" },
    @{ old = "63×97="; new = "12×36=" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
